$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for row 11 (2021年) and row 12 (2022年), columns B..AQ
$row11 = @(9898.389999999999,2788.22,464.78,0.97,3345.98,8626.9,652.62,5249.87,1237.94,1130.87,677.16,397.27,1366.7,4108.04,288.51,1015.74,5181.47,1008.72,14815.38,311.17,4181.14,971.79,11609.08,20311.48,1090.82,2048.45,591.02,2997.58,1842.14,194094.18,35186.45,12029.24,1928.95,1121.35,437.58,8199.9,4946.68,15643.16,542.63,2103.95,3151,593.11)

$row12 = @(11488.5,3028.4,296,1.3,3715.8,9350.299999999999,723.8,6013,1256,1212.2,857.8,455.1,1419.3,4796.4,324.9,1131,5466.7,1312.6,18592.7,367.2,5320.1,1144.3,13502.9,25166.8,1133.9,2548,647.9,3031.6,1839.6,216466.2,36971,12745.2,1991.6,1421.6,524.7,8464.1,3785.3,17158,596.1,2566.9,3454.5,642.9)

# Copy the formatting of the existing last-year label cell (A10) down to the
# two new year-label cells so they keep the bold/centered/bordered style.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 11: 2021年
$ws.Cells.Item(11, 1).Value = "2021年"
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, $i + 2).Value = $row11[$i]
}

# Row 12: 2022年
$ws.Cells.Item(12, 1).Value = "2022年"
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, $i + 2).Value = $row12[$i]
}
